$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B is constant across all new rows: "XmlExtensions.SettingsMenuDef"
$ws.Range("B388:B396").Value = "XmlExtensions.SettingsMenuDef"

# Row 388 - Alchemy Table Drugs?
$ws.Range("C388").Value = "VR_MO_Patch.settings.16.label"
$ws.Range("D388").Value = "Alchemy Table Drugs?"
$ws.Range("E388").Value = "연금술 테이블 약?"
$ws.Range("F388").Value = "SettingsMenuDef: VR_MO_Patch.settings.16.label 'Alchemy Table Drugs?'"

# Row 389 - Allows the alchemy table to inherit drug recipes.
$ws.Range("C389").Value = "VR_MO_Patch.settings.16.tooltip"
$ws.Range("D389").Value = "Allows the alchemy table to inherit drug recipes."
$ws.Range("E389").Value = "연금술 테이블이 약품 레시피를 상속할 수 있도록 허용합니다."
$ws.Range("F389").Value = "SettingsMenuDef: VR_MO_Patch.settings.16.tooltip 'Allows the alchemy table to inherit drug recipes.'"

# Row 390 - Trimmed Meat Options:
$ws.Range("C390").Value = "VR_MO_Patch.settings.18.text"
$ws.Range("D390").Value = "Trimmed Meat Options:"
$ws.Range("E390").Value = "손질된 고기 옵션:"
$ws.Range("F390").Value = "SettingsMenuDef: VR_MO_Patch.settings.18.text 'Trimmed Meat Options:'"

# Row 391 - No trimmed meat.
$ws.Range("C391").Value = "VR_MO_Patch.settings.19.buttons.0.tooltip"
$ws.Range("D391").Value = "No trimmed meat."
$ws.Range("E391").Value = "손질된 고기 없음."
$ws.Range("F391").Value = "SettingsMenuDef: VR_MO_Patch.settings.19.buttons.0.tooltip 'No trimmed meat.'"

# Row 392 - MO Trimmed Meat Recipe:
$ws.Range("C392").Value = "VR_MO_Patch.settings.19.buttons.1.label"
$ws.Range("D392").Value = "MO Trimmed Meat Recipe:"
$ws.Range("E392").Value = "MO 손질된 고기 레시피:"
$ws.Range("F392").Value = "SettingsMenuDef: VR_MO_Patch.settings.19.buttons.1.label 'MO Trimmed Meat Recipe:'"

# Row 393 - This enables a trimmed meat recipe for the extra verisimilitude.
$ws.Range("C393").Value = "VR_MO_Patch.settings.19.buttons.1.tooltip"
$ws.Range("D393").Value = "This enables a trimmed meat recipe for the extra verisimilitude."
$ws.Range("E393").Value = "이렇게 하면 고기를 손질하여 더욱 사실적인 레시피를 만들 수 있습니다."
$ws.Range("F393").Value = "SettingsMenuDef: VR_MO_Patch.settings.19.buttons.1.tooltip 'This enables a trimmed meat recipe for the extra verisimilitude.'"

# Row 394 - Adds various Medieval Overhaul plants to biomes as wild plants.
$ws.Range("C394").Value = "VR_MO_Patch.settings.22.tooltip"
$ws.Range("D394").Value = "Adds various Medieval Overhaul plants to biomes as wild plants."
$ws.Range("E394").Value = "생물군에 다양한 Medieval Overhaul 식물을  야생 식물로 추가합니다."
$ws.Range("F394").Value = "SettingsMenuDef: VR_MO_Patch.settings.22.tooltip 'Adds various Medieval Overhaul plants to biomes as wild plants.'"

# Row 395 - Reset Settings
$ws.Range("C395").Value = "VR_MO_Patch.settings.23.label"
$ws.Range("D395").Value = "Reset Settings"
$ws.Range("E395").Value = "설정 초기화"
$ws.Range("F395").Value = "SettingsMenuDef: VR_MO_Patch.settings.23.label 'Reset Settings'"

# Row 396 - Are you sure?
$ws.Range("C396").Value = "VR_MO_Patch.settings.23.message"
$ws.Range("D396").Value = "Are you sure?"
$ws.Range("E396").Value = "확실하신가요?"
$ws.Range("F396").Value = "SettingsMenuDef: VR_MO_Patch.settings.23.message 'Are you sure?'"

# Column G: pulls the quoted portion out of column F via MID/FIND.
# G388 is a standalone formula; G389:G396 form one shared-formula group.
$ws.Range("G388").Formula = '=MID(F388,FIND("''",F388)+1,FIND("''",F388,FIND("''",F388)+1)-FIND("''",F388)-1)'
$ws.Range("G389:G396").Formula = '=MID(F389,FIND("''",F389)+1,FIND("''",F389,FIND("''",F389)+1)-FIND("''",F389)-1)'

# Restore the cursor/selection to where the author left off.
$ws.Range("D392").Select() | Out-Null
